# Add a "HasPaid" column and five rows of order data, mirroring a
# temp-attachment export that was appended to the order sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column J
$ws.Range("J1").Value = "HasPaid"

# Row data: Id, MemberId, Date, HadDiscount, WasRedeem, Price, CoffeeId, HadAddIn, RedeemId, HasPaid
$rows = @(
    @(1, 3, "1/7/2024", $false, $false, 200, 3, $false, $null, $true),
    @(2, 3, "1/7/2024", $false, $false, 543, 4, $true,  $null, $true),
    @(3, 3, "1/7/2024", $false, $false, 343, 4, $false, $null, $true),
    @(4, 3, "1/7/2024", $false, $false, 343, 4, $false, $null, $false),
    @(5, 3, "1/7/2024", $false, $false, 343, 4, $false, $null, $false)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $r++
}

$ws.Range("A2:J8").Select()

# Touch PageSetup so the worksheet emits a (default, empty) headerFooter
# element, matching a workbook that's been opened/saved in Excel.
$ws.PageSetup.CenterHeader = ""
